$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.413.78'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.16%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.028.55'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.53%  '

$ws.Range('E4').Value = '  +0.20%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '229.25'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.59%  '

$ws.Range('E7').Value = '  +0.02%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '55.78'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.76%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.380'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.15%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0797'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.43%  '

$ws.Range('E11').Value = '  -1.10%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.331.54'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.52%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.34'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.59%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.24'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.34%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.740'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.13%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.20'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.23%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.030.24'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.49%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.347.27'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.28%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.17'
$ws.Range('D19').Style = 'Normal'

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.92'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.27%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0823'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.32%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '223.12'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.16%  '

$ws.Range('E23').Value = '  +0.07%  '

$ws.Range('E24').Value = '  +1.44%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.26'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.19%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.87'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.37%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.11'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.07%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.131'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.82%  '

$ws.Range('E29').Value = '  +0.34%  '

$ws.Range('E30').Value = '  -1.13%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.117'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.68%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.47'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.15%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0604'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.13%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.46'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.00%  '

$ws.Range('E35').Value = '  +8.07%  '

$ws.Range('E36').Value = '  -1.61%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.23'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.69%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.74'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +8.61%  '

$ws.Range('E39').Value = '  +0.07%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.471.47'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.90%  '

$ws.Range('E41').Value = '  -1.43%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0932'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.62%  '

$ws.Range('E43').Value = '  +3.15%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '94.73'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.63%  '

$ws.Range('B45').Value = 'FTXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.22'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +17.05%  '

$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.28'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.87%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.10'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.46%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.00'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.11%  '

$ws.Range('E49').Value = '  +0.84%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.05'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.12%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.220.91'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.41%  '
